# [MOSIP-16645] Preferred languages removed from PreReg UI specs.
#
# The "identity" JSON spec stored in G2 (json_spec column, row 2 = the
# "pre-registration" record) contained a field definition for
# "preferredLang" that needs to be removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("G2")
$jsonSpec = $cell.Value2

$block = "      {`n        `"id`": `"preferredLang`",`n        `"description`": `"Enter your preferred Language`",`n        `"labelName`": {`n          `"eng`": `"Preferred Language`",`n          `"fra`": `"langue préférée`"`n        },`n        `"controlType`": `"dropdown`",`n        `"inputRequired`": true,`n        `"fieldType`": `"dynamic`",`n        `"type`": `"string`",`n        `"validators`": [],`n        `"required`": true`n      },`n"

$hasBlock = $jsonSpec.Contains($block)
if ($hasBlock) {
    $cell.Value2 = $jsonSpec.Replace($block, "")
}

# The two "Fields" title cells (D2 / D5, the identity and fields specs
# rows respectively) pick up a top-vertical-alignment style as part of
# the same formatting pass.
$ws.Range("D2").VerticalAlignment = -4160
$ws.Range("D5").VerticalAlignment = -4160

# Active cell moves from G2 to D2.
[void]$ws.Range("D2").Select()
